# NIT-9010435810.xlsx — "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The previous account-statement ("Estado de Cuenta") detail rows (B16:G30)
# listed arrears for JOSE OSBAIRO GUERRA GONZALEZ, SANDY DAVID BARRIOS MORENO
# and MEYDIS KATHERINE TABORDA PATIÑO across periods 1710-1804, each worker
# appearing once per period in an interleaved order.
#
# This update replaces that block with a refreshed database: the old
# entries are removed and new ones added, now grouped by worker (7
# consecutive periods for SANDY, a single legacy period for JOSE, then 7
# consecutive periods for MEYDIS), with updated "Valor Mora" / "Salario
# Basico" amounts for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico
$data = New-Object "object[,]" 15,6

# SANDY DAVID BARRIOS MORENO - 1047365068 (rows 16-22, periods 1804..1710)
$data[0,0] = "CC"; $data[0,1] = "1047365068"; $data[0,2] = "SANDY DAVID BARRIOS MORENO"; $data[0,3] = "1804"; $data[0,4] = 48000; $data[0,5] = 1200000
$data[1,0] = "CC"; $data[1,1] = "1047365068"; $data[1,2] = "SANDY DAVID BARRIOS MORENO"; $data[1,3] = "1803"; $data[1,4] = 48000; $data[1,5] = 1200000
$data[2,0] = "CC"; $data[2,1] = "1047365068"; $data[2,2] = "SANDY DAVID BARRIOS MORENO"; $data[2,3] = "1802"; $data[2,4] = 48000; $data[2,5] = 1200000
$data[3,0] = "CC"; $data[3,1] = "1047365068"; $data[3,2] = "SANDY DAVID BARRIOS MORENO"; $data[3,3] = "1801"; $data[3,4] = 48000; $data[3,5] = 1200000
$data[4,0] = "CC"; $data[4,1] = "1047365068"; $data[4,2] = "SANDY DAVID BARRIOS MORENO"; $data[4,3] = "1712"; $data[4,4] = 48000; $data[4,5] = 1200000
$data[5,0] = "CC"; $data[5,1] = "1047365068"; $data[5,2] = "SANDY DAVID BARRIOS MORENO"; $data[5,3] = "1711"; $data[5,4] = 48000; $data[5,5] = 1200000
$data[6,0] = "CC"; $data[6,1] = "1047365068"; $data[6,2] = "SANDY DAVID BARRIOS MORENO"; $data[6,3] = "1710"; $data[6,4] = 48000; $data[6,5] = 1200000

# JOSE OSBAIRO GUERRA GONZALEZ - 71184413 (row 23, period 1710)
$data[7,0] = "CC"; $data[7,1] = "71184413"; $data[7,2] = "JOSE OSBAIRO GUERRA GONZALEZ"; $data[7,3] = "1710"; $data[7,4] = 48000; $data[7,5] = 1200000

# MEYDIS KATHERINE TABORDA PATIÑO - 1047451676 (rows 24-30, periods 1804..1710)
$data[8,0]  = "CC"; $data[8,1]  = "1047451676"; $data[8,2]  = "MEYDIS KATHERINE TABORDA PATIÑO"; $data[8,3]  = "1804"; $data[8,4]  = 64000; $data[8,5]  = 1600000
$data[9,0]  = "CC"; $data[9,1]  = "1047451676"; $data[9,2]  = "MEYDIS KATHERINE TABORDA PATIÑO"; $data[9,3]  = "1803"; $data[9,4]  = 64000; $data[9,5]  = 1600000
$data[10,0] = "CC"; $data[10,1] = "1047451676"; $data[10,2] = "MEYDIS KATHERINE TABORDA PATIÑO"; $data[10,3] = "1802"; $data[10,4] = 64000; $data[10,5] = 1600000
$data[11,0] = "CC"; $data[11,1] = "1047451676"; $data[11,2] = "MEYDIS KATHERINE TABORDA PATIÑO"; $data[11,3] = "1801"; $data[11,4] = 64000; $data[11,5] = 1600000
$data[12,0] = "CC"; $data[12,1] = "1047451676"; $data[12,2] = "MEYDIS KATHERINE TABORDA PATIÑO"; $data[12,3] = "1712"; $data[12,4] = 64000; $data[12,5] = 1600000
$data[13,0] = "CC"; $data[13,1] = "1047451676"; $data[13,2] = "MEYDIS KATHERINE TABORDA PATIÑO"; $data[13,3] = "1711"; $data[13,4] = 64000; $data[13,5] = 1600000
$data[14,0] = "CC"; $data[14,1] = "1047451676"; $data[14,2] = "MEYDIS KATHERINE TABORDA PATIÑO"; $data[14,3] = "1710"; $data[14,4] = 55466; $data[14,5] = 1600000

$ws.Range("B16:G30").Value = $data
